$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 8999.5
$ws.Range("I20").Value = 8999.5
$ws.Range("K20").Value = 8999.5
$ws.Range("M20").Value = -8769.5

$ws.Range("H35").Value = 8999.5
$ws.Range("I35").Value = 8999.5
$ws.Range("K35").Value = 8999.5
$ws.Range("M35").Value = -8620.5

$ws.Range("H137").Value = 3190.8462
$ws.Range("I137").Value = 2686.2
$ws.Range("J137").Value = 3506.25
$ws.Range("K137").Value = 8058.599999999999
$ws.Range("L137").Value = 10518.75
$ws.Range("M137").Value = -5508.599999999999
$ws.Range("N137").Value = -15618.75

$ws.Range("H138").Value = 2308.85
$ws.Range("I138").Value = 930.02325
$ws.Range("J138").Value = 3349.0176
$ws.Range("K138").Value = 2790.06975
$ws.Range("L138").Value = 10047.0528
$ws.Range("M138").Value = 2349.93025
$ws.Range("N138").Value = -20327.0528

$ws.Range("H141").Value = 1921.55
$ws.Range("I141").Value = 1525.3529
$ws.Range("K141").Value = 4576.0587
$ws.Range("M141").Value = 603.9412999999995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 8298.8
$ws.Range("I2").Value = 973.3
$ws.Range("K2").Value = 973.3
$ws.Range("M2").Value = -860.3

$ws.Range("H32").Value = 2228.762
$ws.Range("I32").Value = 2119.9507
$ws.Range("K32").Value = 2119.9507
$ws.Range("M32").Value = -1832.9507

$ws.Range("H61").Value = 5019.5835
$ws.Range("I61").Value = 2459.8096
$ws.Range("J61").Value = 22938
$ws.Range("K61").Value = 2459.8096
$ws.Range("L61").Value = 22938
$ws.Range("M61").Value = -2247.8096
$ws.Range("N61").Value = -23362

$ws.Range("H74").Value = 14496202
$ws.Range("I74").Value = 15875494
$ws.Range("J74").Value = 13632
$ws.Range("K74").Value = 15875494
$ws.Range("L74").Value = 13632
$ws.Range("M74").Value = -15874620
$ws.Range("N74").Value = -15380

$ws.Range("H77").Value = 14496202
$ws.Range("I77").Value = 15875494
$ws.Range("J77").Value = 13632
$ws.Range("K77").Value = 79377470
$ws.Range("L77").Value = 68160
$ws.Range("M77").Value = -79373102
$ws.Range("N77").Value = -76896

$ws.Range("H116").Value = 8298.8
$ws.Range("I116").Value = 973.3
$ws.Range("K116").Value = 973.3
$ws.Range("M116").Value = 1320.7

$ws.Range("H132").Value = 4181.091
$ws.Range("I132").Value = 3483.1516
$ws.Range("K132").Value = 10449.4548
$ws.Range("M132").Value = -7919.4548

$ws.Range("H136").Value = 5019.5835
$ws.Range("I136").Value = 2459.8096
$ws.Range("J136").Value = 22938
$ws.Range("K136").Value = 7379.4288
$ws.Range("L136").Value = 68814
$ws.Range("M136").Value = -4829.4288
$ws.Range("N136").Value = -73914

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 8298.8
$ws.Range("I3").Value = 973.3
$ws.Range("K3").Value = 973.3
$ws.Range("M3").Value = -859.3

$ws.Range("H134").Value = 1831
$ws.Range("I134").Value = 1831
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 5493
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -2958
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2219.0217
$ws.Range("I58").Value = 1035.8684
$ws.Range("J58").Value = 7839
$ws.Range("K58").Value = 1035.8684
$ws.Range("L58").Value = 7839
$ws.Range("M58").Value = -832.8684000000001
$ws.Range("N58").Value = -8245

$ws.Range("H134").Value = 1800.4648
$ws.Range("I134").Value = 1145.4584
$ws.Range("J134").Value = 3167.4348
$ws.Range("K134").Value = 3436.3752
$ws.Range("L134").Value = 9502.3044
$ws.Range("M134").Value = -901.3752
$ws.Range("N134").Value = -14572.3044

$ws.Range("H136").Value = 2219.0217
$ws.Range("I136").Value = 1035.8684
$ws.Range("J136").Value = 7839
$ws.Range("K136").Value = 3107.6052
$ws.Range("L136").Value = 23517
$ws.Range("M136").Value = -557.6052
$ws.Range("N136").Value = -28617

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1077476.8
$ws.Range("I5").Value = 1398.4667
$ws.Range("J5").Value = 2086300.2
$ws.Range("K5").Value = 4195.4001
$ws.Range("L5").Value = 6258900.6
$ws.Range("M5").Value = -4083.4001
$ws.Range("N5").Value = -6259124.6

$ws.Range("H12").Value = 65.47059
$ws.Range("I12").Value = 148.6
$ws.Range("J12").Value = 30.833334
$ws.Range("K12").Value = 445.8
$ws.Range("L12").Value = 92.500002
$ws.Range("M12").Value = -272.8
$ws.Range("N12").Value = -438.500002

$ws.Range("H60").Value = 15525439
$ws.Range("I60").Value = 23809672
$ws.Range("J60").Value = 1028029.4
$ws.Range("K60").Value = 71429016
$ws.Range("L60").Value = 3084088.2
$ws.Range("M60").Value = -71428765
$ws.Range("N60").Value = -3084590.2

$ws.Range("H135").Value = 1077476.8
$ws.Range("I135").Value = 1398.4667
$ws.Range("J135").Value = 2086300.2
$ws.Range("K135").Value = 12586.2003
$ws.Range("L135").Value = 18776701.8
$ws.Range("M135").Value = -10051.2003
$ws.Range("N135").Value = -18781771.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1441.1428
$ws.Range("I97").Value = 1507.3334
$ws.Range("J97").Value = 1242.5714
$ws.Range("K97").Value = 1507.3334
$ws.Range("L97").Value = 1242.5714
$ws.Range("M97").Value = -1011.3334
$ws.Range("N97").Value = -2234.5714

$ws.Range("H132").Value = 3726.653
$ws.Range("I132").Value = 3435.122
$ws.Range("K132").Value = 10305.366
$ws.Range("M132").Value = -7775.366

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1509.4348
$ws.Range("I16").Value = 1509.4348
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1509.4348
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1339.4348
$ws.Range("N16").ClearContents()

$ws.Range("H132").Value = 4753.3335
$ws.Range("I132").Value = 4087.375
$ws.Range("J132").Value = 10081
$ws.Range("K132").Value = 12262.125
$ws.Range("L132").Value = 30243
$ws.Range("M132").Value = -9732.125
$ws.Range("N132").Value = -35303

$ws.Range("H136").Value = 4048.5854
$ws.Range("I136").Value = 1916.3939
$ws.Range("J136").Value = 12843.875
$ws.Range("K136").Value = 5749.1817
$ws.Range("L136").Value = 38531.625
$ws.Range("M136").Value = -3199.1817
$ws.Range("N136").Value = -43631.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1518.9773
$ws.Range("I136").Value = 1130.5585
$ws.Range("J136").Value = 4237.909
$ws.Range("K136").Value = 3391.6755
$ws.Range("L136").Value = 12713.727
$ws.Range("M136").Value = -841.6755000000003
$ws.Range("N136").Value = -17813.727
